$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A10 ("2024-01-03") looks like a date, so Excel would otherwise coerce it to
# a date serial number. Briefly mark the cell as Text, assign the literal
# string, then clear the formatting again so no stray style index is left
# behind on the cell (matches rows 2-9, which carry no "s" attribute).
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "2024-01-03"
$ws.Range("A10").ClearFormats()

$ws.Range("B10").Value = "09:18:46"
$ws.Range("C10").Value = "Wednesday"

# D10 ("00") looks like a plain number, so it needs the same treatment or it
# would be stored as 0 instead of the text "00".
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "00"
$ws.Range("D10").ClearFormats()

$ws.Range("E10").Value = 140196
$ws.Range("F10").Value = 142836
$ws.Range("G10").Value = 171191
$ws.Range("H10").Value = 146039
$ws.Range("I10").Value = -1
$ws.Range("J10").Value = 116840
$ws.Range("K10").Value = 223604
$ws.Range("L10").Value = 247829
$ws.Range("M10").Value = 183765
$ws.Range("N10").Value = 109777
$ws.Range("O10").Value = 39740
$ws.Range("P10").Value = 30775
$ws.Range("Q10").Value = 71944
$ws.Range("R10").Value = -1
$ws.Range("S10").Value = 40848
$ws.Range("T10").Value = -1
